$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range('A1').Value = 'Datos actualizados a 1 de Abril de 2020 a las 13:50'

# Province rows that moved up/down in the "Casos totales" sort order
# and/or received refreshed case counts. Each entry is
# Row -> (Ciudad, Casos totales, Casos activos, Recuperados, Muertes)
$updates = @(
    @{ Row = 9; Ciudad = 'Navarra'; Totales = 2497; Activos = 236; Recuperados = 2131; Muertes = 130 }
    @{ Row = 10; Ciudad = 'Ciudad Real'; Totales = 2471; Activos = 397; Recuperados = 2076; Muertes = 245 }
    @{ Row = 18; Ciudad = 'Malaga'; Totales = 1464; Activos = 93; Recuperados = 1294; Muertes = 77 }
    @{ Row = 19; Ciudad = 'Pontevedra'; Totales = 1452; Activos = 259; Recuperados = 1337; Muertes = 30 }
    @{ Row = 29; Ciudad = 'Sevilla'; Totales = 1215; Activos = 20; Recuperados = 1140; Muertes = 55 }
    @{ Row = 30; Ciudad = 'Cantabria'; Totales = 1213; Activos = 43; Recuperados = 1116; Muertes = 54 }
    @{ Row = 31; Ciudad = 'Gipuzkoa/Guipuzcoa'; Totales = 1206; Activos = 2165; Recuperados = 639; Muertes = 52 }
    @{ Row = 32; Ciudad = 'Granada'; Totales = 1182; Activos = 15; Recuperados = 1087; Muertes = 80 }
    @{ Row = 33; Ciudad = 'Valladolid'; Totales = 1109; Activos = 262; Recuperados = 758; Muertes = 89 }
    @{ Row = 34; Ciudad = 'Caceres'; Totales = 1093; Activos = 31; Recuperados = 932; Muertes = 130 }
    @{ Row = 41; Ciudad = 'Jaen'; Totales = 742; Activos = 17; Recuperados = 688; Muertes = 37 }
    @{ Row = 42; Ciudad = 'Cordoba'; Totales = 687; Activos = 4; Recuperados = 662; Muertes = 21 }
    @{ Row = 45; Ciudad = 'Cadiz'; Totales = 645; Activos = 17; Recuperados = 611; Muertes = 17 }
    @{ Row = 46; Ciudad = 'Ourense'; Totales = 626; Activos = 259; Recuperados = 559; Muertes = 18 }
    @{ Row = 47; Ciudad = 'Badajoz'; Totales = 586; Activos = 82; Recuperados = 482; Muertes = 22 }
    @{ Row = 54; Ciudad = 'Almeria'; Totales = 262; Activos = 14; Recuperados = 231; Muertes = 17 }
    @{ Row = 57; Ciudad = 'Huelva'; Totales = 195; Activos = 2; Recuperados = 189; Muertes = 4 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Ciudad
    $ws.Cells.Item($r, 2).Value = $u.Totales
    $ws.Cells.Item($r, 3).Value = $u.Activos
    $ws.Cells.Item($r, 4).Value = $u.Recuperados
    $ws.Cells.Item($r, 5).Value = $u.Muertes
}
